$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template Modelo")

# Insert a new column before column E ("BOX" transaction column), shifting
# REFERENCE_ID's neighbours right. The new column inherits formatting from
# column D, matching Excel's default column-insert behaviour.
$colDWidth = $ws.Columns("D:D").ColumnWidth
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = $colDWidth

# Header label for the newly inserted column.
$ws.Range("E1").Value = "BOX"

# Update the selection to mirror the post-edit cursor position.
$ws.Range("E2").Select()
